$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table by one year column: copy the formatting of column P
# (2019) into the new column Q before writing the 2020 values, so the new
# cells inherit the same number formats / fonts / borders as the rest of
# the table.
$ws.Range("P4:P14").Copy()
$ws.Range("Q4:Q14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New "2020" column of data.
$ws.Range("Q4").Value = 2020
$ws.Range("Q5").Value = 109.7221295941265
$ws.Range("Q6").Value = 108.44905375816947
$ws.Range("Q7").Value = 109.90982951756889
$ws.Range("Q8").Value = 108.40606487500015
$ws.Range("Q9").Value = 109.40161876466024
$ws.Range("Q10").Value = 107.71155656686271
$ws.Range("Q11").Value = 111.78921596090774
$ws.Range("Q12").Value = 111.39254046803097
$ws.Range("Q13").Value = 110.44919152842827
$ws.Range("Q14").Value = 106.89826464456031

# Match the author's final selection (cell N14 active) when the file was
# last saved.
[void]$ws.Range("N14").Select()

Write-Host "edit applied"
